$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.860.30"
$ws.Range("E2").Value = "  +1.33%  "
$ws.Range("D3").Value = "3.649.61"
$ws.Range("E3").Value = "  +6.29%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'581.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'176.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.38%  "
$ws.Range("D7").Value = "3.640.76"
$ws.Range("E7").Value = "  +6.25%  "
$ws.Range("D8").Value = "'0.619"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("D10").Value = "'0.199"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.91%  "
$ws.Range("D11").Value = "'6.90"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +24.69%  "
$ws.Range("D12").Value = "'0.607"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.74%  "
$ws.Range("D13").Value = "'48.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("E14").Value = "  +2.22%  "
$ws.Range("D15").Value = "4.235.37"
$ws.Range("E15").Value = "  +6.32%  "
$ws.Range("D16").Value = "'674.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").Value = "'8.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.54%  "
$ws.Range("D18").Value = "3.660.97"
$ws.Range("E18").Value = "  +6.68%  "
$ws.Range("D19").Value = "70.893.82"
$ws.Range("E19").Value = "  +1.39%  "
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("D21").Value = "'17.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'11.51"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "
$ws.Range("E23").Value = "  +4.01%  "
$ws.Range("D24").Value = "'17.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.42%  "
$ws.Range("D25").Value = "'101.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").Value = "'2.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.43%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "'34.96"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.10%  "
$ws.Range("E31").Value = "  +3.65%  "
$ws.Range("E32").Value = "  +3.43%  "
$ws.Range("E33").Value = "  -3.06%  "
$ws.Range("D34").Value = "'7.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.76%  "
$ws.Range("D35").Value = "'4.01"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.24%  "
$ws.Range("D36").Value = "'581.58"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.03%  "
$ws.Range("D37").Value = "'11.13"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  +3.96%  "
$ws.Range("D39").Value = "'58.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").Value = "3.605.95"
$ws.Range("E41").Value = "  +0.37%  "
$ws.Range("D42").Value = "'0.0457"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.65%  "
$ws.Range("E43").Value = "  +1.88%  "
$ws.Range("E44").Value = "  +4.47%  "
$ws.Range("D45").Value = "0.0₃0749"
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("D46").Value = "'34.97"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").Value = "'2.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").Value = "'2.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.37%  "
$ws.Range("D49").Value = "'0.134"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.57%  "
$ws.Range("D50").Value = "'135.10"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("D51").Value = "'2.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +9.45%  "
